# sim-explorer.pptx edit:
#  - refresh the cached "DocumentDate" field text (11 November 2024 -> 19
#    November 2024) everywhere it is placed: handout master, notes master,
#    and the three slide layouts ("Three Content", "Three Content, grey",
#    ">Do not use layouts after this >") that carry their own copy of the
#    placeholder.
#  - tweak the "Set and get variables" slide (slide 7) syntax blurb so the
#    right-hand side reads "<value(s)>" instead of "<value>".

$p = $ppt.ActivePresentation

$oldDate = "11 November 2024"
$newDate = "19 November 2024"

# --- Handout master: "Date Placeholder 8" ---------------------------------
$hm = $p.HandoutMaster
for ($i = 1; $i -le $hm.Shapes.Count; $i++) {
    $sh = $hm.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

# --- Notes master: "Date Placeholder 8" -----------------------------------
$nm = $p.NotesMaster
for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
    $sh = $nm.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

# --- Slide-master custom layouts: only the real date *field* shapes ------
# (Layout 17 also has a static, user-drawn "SD_FLD_DocumentDate" textbox
# with the same literal text - that one is left untouched, matching the
# diff, so match by placeholder name rather than by raw text.)
$sm = $p.SlideMaster
foreach ($layoutIdx in 16, 17, 29) {
    $cl = $sm.CustomLayouts.Item($layoutIdx)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $sh = $cl.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.Name -like "Date Placeholder*" -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Slide 7 ("Set and get variables"): <value> -> <value(s)> ------------
$s7 = $p.Slides.Item(7)
for ($i = 1; $i -le $s7.Shapes.Count; $i++) {
    $sh = $s7.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "<variable-name><slice><@time> : <value>") {
        $sh.TextFrame.TextRange.Text = "<variable-name><slice><@time> : <value(s)>"
    }
}
